{"js": "// Insert a centered contact-info paragraph right after the \"Dheeraj Chand\"\n// name/title paragraph (and before \"PROFESSIONAL SUMMARY\"), matching the\n// look of a plain (non-bold, default-size) run the way the author's other\n// generated resumes render their contact line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet namePara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Dheeraj Chand\") {\n    namePara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!namePara) {\n  throw new Error('Could not locate the \"Dheeraj Chand\" paragraph.');\n}\n\nconst contactText =\n  \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\";\n\n// Create a fresh (initially empty) paragraph right after the name line.\n// insertParagraph() would normally copy the name line's bold/large-size run\n// formatting, so instead we overwrite the brand-new paragraph's OOXML\n// outright via insertOoxml(\"Replace\") -- that gives us a clean paragraph\n// with only centered alignment and a plain run, exactly like the rest of\n// the document's body text.\nconst newPara = namePara.insertParagraph(\"\", \"After\");\nconst newRange = newPara.getRange();\n\nconst escapeXml = (s) =>\n  s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n\nconst ooxml = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr><w:jc w:val=\"center\"/></w:pPr>\n            <w:r><w:t>${escapeXml(contactText)}</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nnewRange.insertOoxml(ooxml, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Insert a centered contact-info paragraph right after the \"Dheeraj Chand\"\n# name/title paragraph (and before \"PROFESSIONAL SUMMARY\"), matching the\n# look of a plain (non-bold, default-size) run like the rest of the body\n# text uses.\n$d = $word.ActiveDocument\n\n$namePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Dheeraj Chand\") {\n        $namePara = $p\n        break\n    }\n}\n\nif ($null -eq $namePara) {\n    throw \"Could not locate the 'Dheeraj Chand' paragraph.\"\n}\n\n$contactText = \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\"\n\n# Mint a brand-new (empty) paragraph right after the name line. A plain\n# $newRange.Text assignment would inherit the name line's bold/28pt run\n# formatting, so instead we overwrite the fresh paragraph's underlying OOXML\n# outright via Range.InsertXML -- that produces a clean paragraph with only\n# centered alignment and a plain run, matching the rest of the document's\n# body text.\n$namePara.Range.InsertParagraphAfter()\n$newPara = $namePara.Next()\n$newRange = $newPara.Range\n\n$ooxml = @\"\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr><w:jc w:val=\"center\"/></w:pPr>\n            <w:r><w:t>$contactText</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$newRange.InsertXML($ooxml)\n"}
